# Apply the "latest data" refresh to the eventi_modena sheet:
#  1. Row 10 gets a refreshed last-modified timestamp (L10) and a reworded
#     schedule blurb (N10) for the "Kassandra" event - no row movement.
#  2. A brand-new event row ("HEREAFTER. Episodi teatrali") is inserted at
#     row 16, pushing every existing row from 16 down to 51 one row lower
#     (17 down to 52). The sheet's used range grows from A1:AA51 to A1:AA52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. In-place edits on row 10 -------------------------------------------
$ws.Range("L10").Value = "2022-06-06T07:38:50+00:00"
$ws.Range("N10").Value = " 31 maggio e 1 giugno ore 20.30  repliche dal 2 al 12 giugno ore 21.15"

# --- 2. Insert a new row at 16 and shift everything below it down ---------
$ws.Rows("16:16").Insert()

# --- 3. Populate the newly inserted row 16 with the new event's data ------
$ws.Range("A16").Value  = "Spettacoli,Teatro"
$ws.Range("B16").Value  = "Modena"
$ws.Range("C16").Value  = "Viale Caduti in Guerra, 196"
$ws.Range("D16").Value  = "2020-09-17T12:45:54+00:00"
$ws.Range("E16").Value  = "A cura delle allieve attrici e degli allievi attori di ERT / Teatro Nazionale"
$ws.Range("F16").Value  = "2014-09-30T12:50:00+00:00"
$ws.Range("G16").Value  = "info@emiliaromagnateatro.com"
$ws.Range("H16").Value  = "2022-06-07T22:00:00+00:00"
$ws.Range("I16").Value  = "2022-06-11T21:55:00+00:00"
$ws.Range("J16").Value  = "https://www.comune.modena.it/api/novita/eventi/2022/hereafter-episodi-teatrali/@@images/59c6e256-ecc2-48df-912f-b171b831bd83.jpeg"
$ws.Range("K16").Value  = ""
$ws.Range("L16").Value  = "2022-06-06T07:48:28+00:00"
$ws.Range("M16").Value  = "Teatro Tempio"
$ws.Range("N16").Value  = " ore 19.00"
$ws.Range("O16").Value  = ""
$ws.Range("P16").Value  = " A pagamento, vedi nel testo le info per l'acquisto dei biglietti."
$ws.Range("Q16").Value  = ""
$ws.Range("R16").Value  = "059/2163021"
$ws.Range("S16").Value  = "HEREAFTER. Episodi teatrali"
$ws.Range("T16").Value  = ""
$ws.Range("U16").Value  = "http://www.emiliaromagnateatro.com"
$ws.Range("V16").Value  = $false
$ws.Range("W16").Value  = ""
$ws.Range("X16").Value  = "https://www.comune.modena.it/novita/eventi/2022/hereafter-episodi-teatrali"
$ws.Range("Y16").Value  = "44,64381951149482"
$ws.Range("Z16").Value  = "10,93139345085676"
$ws.Range("AA16").Value = "POINT (10.93139345085676 44.64381951149482)"
